# Camara 1.1 (Foto en BBDD)
#
# The "Planificador" task sheet gains two new tracked tasks, inserted
# right above the "Ruleta Europea" task (previously row 47):
#   - "Poner la foto en la BBDD"      (Alvaro, Nº Historia 31, 3h reales)
#   - "Interfaz Funciones Admin"      (Alvaro, no hours logged yet)
#
# Inserting the two rows natively (EntireRow.Insert) makes Excel itself
# re-point every formula/shared-formula range that lived below the
# insertion point (K3, K5, K6, K7, K8, the E3:F3 shared SUM) — exactly
# like a human using "Insert Rows" in the UI, so we don't hand-edit any
# formula text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planificador")
$ws.Activate()

# --- Insert two blank rows above the old row 47 ("Ruleta Europea") ---
$ws.Range("A47:A48").EntireRow.Insert()

# Pick up the same visual style (fonts/fills/borders) the surrounding
# task rows use: column B uses style "8" (owner name), C:G use style "1".
# The cleanest way to reproduce that exactly is to copy formats from the
# row immediately below (old row 47, now shifted down to row 49, which
# still carries that same style pattern).
$ws.Range("B49:G49").Copy()
$ws.Range("B47:G47").PasteSpecial(-4122)
$ws.Range("B49:G49").Copy()
$ws.Range("B48:G48").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new task data ---
# Row 48's text is entered first so the shared-string table gets the
# same ordering as the authored workbook (Interfaz Funciones Admin,
# then Poner la foto en la BBDD).
$ws.Cells.Item(48, 2).Value2 = "Álvaro"
$ws.Cells.Item(48, 3).Value2 = "Interfaz Funciones Admin"

$ws.Cells.Item(47, 2).Value2 = "Álvaro"
$ws.Cells.Item(47, 3).Value2 = "Poner la foto en la BBDD"
$ws.Cells.Item(47, 4).Value2 = 31
$ws.Cells.Item(47, 6).Value2 = 3

# --- Reflect where the author ended up looking on the sheet ---
$ws.Range("I41").Select()
$excel.ActiveWindow.Zoom = 85
